# Update column G ("K") values on the active sheet to reflect the
# regenerated strikeout ("K") counts instead of the old "Strike#" values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 2
    4  = 3
    5  = 2
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
